# Update the "Estado de Cuenta" worksheet:
#  - Column E (Periodo Mora) for rows 16-24 is re-sorted ascending
#    (2210, 2211, 2212, 2301, 2302, 2303, 2304, 2305, 2306) instead of the
#    previous descending order.
#  - Column G (Salario Basico) for rows 16-24 is updated from 1300000 to
#    1423500 (new base salary amount for this batch of records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2210", "2211", "2212", "2301", "2302", "2303", "2304", "2305", "2306")
$salario = 1423500

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("G$row").Value = $salario
}
